# Insert a new weekly record as row 40, pushing the existing rows 40-65
# down to 41-66 (a new "Fruta, Vega Modelo de Temuco - Caqui" reading).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a fresh row at position 40.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly reading.
$ws.Cells.Item(40, 1).Value = 10
$ws.Cells.Item(40, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(40, 3).Value = "La Araucanía"
$ws.Cells.Item(40, 4).Value = 45062
$ws.Cells.Item(40, 5).Value = 9
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100107
$ws.Cells.Item(40, 8).Value = "Otros"
$ws.Cells.Item(40, 9).Value = 100107001
$ws.Cells.Item(40, 10).Value = "Caqui"
$ws.Cells.Item(40, 11).Value = "Fuyu"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 100
$ws.Cells.Item(40, 14).Value = 18000
$ws.Cells.Item(40, 15).Value = 18000
$ws.Cells.Item(40, 16).Value = 18000
$ws.Cells.Item(40, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(40, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(40, 19).Value = 1200
$ws.Cells.Item(40, 20).Value = 15
